# Updated data and place order
#
# 1. "Sheet1": the source-of-truth counter cell I23 bumps from 5 -> 6.
#    Every name/email in A23:C42 is produced by CONCATENATE(...,$I$23[,...])
#    so they all ripple automatically once I23 changes.
# 2. "order": columns R/S/T (rows 2-21) hold literal copies of the same
#    "<Name><n>" / "<Name><n>@gmail.com" strings (not formulas), so they
#    need to be rewritten explicitly with the new suffix.
# 3. "order" sheetView selection moves from R5 to R2:T21 (active cell R2).

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$order  = $wb.Worksheets.Item("order")

$names = @(
    "DonnellJernigan",
    "MalikOtoole",
    "AlanCaudill",
    "AdanApplegate",
    "AiyanaWhitworth",
    "MercedezBrien",
    "DuaneHager",
    "LorenBell",
    "GeraldHiller",
    "DeionBranch",
    "DakotaHalstead",
    "ElliottFurman",
    "MiltonCamp",
    "DawnChester",
    "ZacheryPetrie",
    "EstebanAngel",
    "JimmyBlankenship",
    "AllysaGrice",
    "AugustineYoo",
    "BrandiSouthard"
)

# Bump the counter that drives Sheet1!A23:C42 through their CONCATENATE formulas.
$sheet1.Range("I23").Value = 6

# Rewrite the literal name/email copies on the "order" sheet (R2:T21).
for ($i = 0; $i -lt $names.Count; $i++) {
    $row = 2 + $i
    $newName = $names[$i] + "6"
    $newEmail = $newName + "@gmail.com"

    $order.Cells.Item($row, 18).Value = $newName   # R
    $order.Cells.Item($row, 19).Value = $newName   # S
    $order.Cells.Item($row, 20).Value = $newEmail  # T
}

# Move the active selection on the "order" sheet to R2:T21.
$order.Activate()
$order.Range("R2:T21").Select()
